$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: update row 3 (d3775bf9... file) handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 20:37:51"
$wsZhCn.Range("H3").Value = "2016-03-19 20:38:10"

# "de-de" worksheet: update row 3 (d3775bf9... file) handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 20:37:53"
$wsDeDe.Range("H3").Value = "2016-03-19 20:38:15"
